$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.646.43'
$ws.Range('E2').Value = '  -3.20%  '
$ws.Range('D3').Value = '1.739.10'
$ws.Range('E3').Value = '  -5.57%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.67'
$ws.Range('E5').Value = '  -10.23%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4896'
$ws.Range('E7').Value = '  -7.99%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.23'
$ws.Range('E8').Value = '  -8.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2544'
$ws.Range('E9').Value = '  -17.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06081'
$ws.Range('E10').Value = '  -11.85%  '
$ws.Range('D11').Value = '1.739.30'
$ws.Range('E11').Value = '  -5.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06844'
$ws.Range('E12').Value = '  -12.54%  '
$ws.Range('E13').Value = '  -20.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.430'
$ws.Range('E15').Value = '  -15.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.5580'
$ws.Range('E16').Value = '  -26.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').Value = '25.680.63'
$ws.Range('E19').Value = '  -3.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.40'
$ws.Range('E20').Value = '  -18.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006564'
$ws.Range('E21').Value = '  -17.37%  '
$ws.Range('D22').Value = '1.961.47'
$ws.Range('E22').Value = '  -6.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.026'
$ws.Range('E23').Value = '  -13.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.892'
$ws.Range('E24').Value = '  -15.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.990'
$ws.Range('E25').Value = '  -16.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '137.23'
$ws.Range('E26').Value = '  -3.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.479'
$ws.Range('E27').Value = '  -12.30%  '
$ws.Range('E28').Value = '  -16.90%  '
$ws.Range('E29').Value = '  -13.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '100.99'
$ws.Range('E30').Value = '  -9.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07975'
$ws.Range('E31').Value = '  -9.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.671'
$ws.Range('E32').Value = '  -14.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.364'
$ws.Range('E33').Value = '  -17.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04401'
$ws.Range('E34').Value = '  -8.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.001'
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.612'
$ws.Range('E36').Value = '  -10.97%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9646'
$ws.Range('E37').Value = '  -14.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5862'
$ws.Range('E38').Value = '  -19.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.636'
$ws.Range('E39').Value = '  -14.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.002'
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '103.49'
$ws.Range('E41').Value = '  -4.43%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01497'
$ws.Range('E42').Value = '  -13.18%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.857'
$ws.Range('E43').Value = '  -19.98%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.125'
$ws.Range('E44').Value = '  -12.96%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3708'
$ws.Range('E45').Value = '  -22.72%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.7205'
$ws.Range('E46').Value = '  -20.20%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05208'
$ws.Range('E47').Value = '  -10.26%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1084'
$ws.Range('E48').Value = '  -12.67%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '29.83'
$ws.Range('E49').Value = '  -14.68%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '51.59'
$ws.Range('E50').Value = '  -14.44%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.752'
$ws.Range('E51').Value = '  -24.57%  '
